# Auto-generated edit script
# Applies cached market-data value updates (scheduled runner refresh)
# to the leve-profit tracking workbook, matching the upstream commit's
# per-cell diff across all 8 job sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 12084.267
$ws.Range("I51").Value = 2955
$ws.Range("J51").Value = 12736.357
$ws.Range("K51").Value = 2955
$ws.Range("L51").Value = 12736.357
$ws.Range("M51").Value = -2471
$ws.Range("N51").Value = -13704.357
$ws.Range("H53").Value = 228
$ws.Range("I53").Value = 149.6923
$ws.Range("K53").Value = 149.6923
$ws.Range("M53").Value = 487.3077
$ws.Range("H129").Value = 2382.389
$ws.Range("I129").Value = 799
$ws.Range("J129").Value = 2699.0667
$ws.Range("K129").Value = 2397
$ws.Range("L129").Value = 8097.2001
$ws.Range("M129").Value = 2603
$ws.Range("N129").Value = -18097.2001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 397.33334
$ws.Range("I2").Value = 397
$ws.Range("K2").Value = 397
$ws.Range("M2").Value = -284
$ws.Range("H32").Value = 9976.866
$ws.Range("I32").Value = 11242.454
$ws.Range("J32").Value = 6496.5
$ws.Range("K32").Value = 11242.454
$ws.Range("L32").Value = 6496.5
$ws.Range("M32").Value = -10955.454
$ws.Range("N32").Value = -7070.5
$ws.Range("H45").Value = 4487.1113
$ws.Range("I45").Value = 3796.25
$ws.Range("K45").Value = 3796.25
$ws.Range("M45").Value = -3419.25
$ws.Range("H88").Value = 687.25
$ws.Range("I88").Value = 624.5
$ws.Range("J88").Value = 750
$ws.Range("K88").Value = 624.5
$ws.Range("L88").Value = 750
$ws.Range("M88").Value = -218.5
$ws.Range("N88").Value = -1562
$ws.Range("H91").Value = 687.25
$ws.Range("I91").Value = 624.5
$ws.Range("J91").Value = 750
$ws.Range("K91").Value = 624.5
$ws.Range("L91").Value = 750
$ws.Range("M91").Value = 779.5
$ws.Range("N91").Value = -3558
$ws.Range("H97").Value = 309
$ws.Range("I97").Value = 309
$ws.Range("K97").Value = 309
$ws.Range("M97").Value = 187
$ws.Range("H110").Value = 2500
$ws.Range("I110").Value = 2500
$ws.Range("K110").Value = 2500
$ws.Range("M110").Value = -455
$ws.Range("H116").Value = 397.33334
$ws.Range("I116").Value = 397
$ws.Range("K116").Value = 397
$ws.Range("M116").Value = 1897
$ws.Range("H122").Value = 913
$ws.Range("I122").Value = 812
$ws.Range("J122").Value = 1014
$ws.Range("K122").Value = 2436
$ws.Range("L122").Value = 3042
$ws.Range("M122").Value = 14
$ws.Range("N122").Value = -7942

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 397.33334
$ws.Range("I3").Value = 397
$ws.Range("K3").Value = 397
$ws.Range("M3").Value = -283
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H86").Value = 5435
$ws.Range("I86").Value = 599.5
$ws.Range("K86").Value = 599.5
$ws.Range("M86").Value = 523.5
$ws.Range("H89").Value = 5435
$ws.Range("I89").Value = 599.5
$ws.Range("K89").Value = 2997.5
$ws.Range("M89").Value = 2618.5
$ws.Range("H105").Value = 3999
$ws.Range("I105").Value = 3999
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3999
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -2252
$ws.Range("N105").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3758.1738
$ws.Range("I105").Value = 3040.7
$ws.Range("K105").Value = 3040.7
$ws.Range("M105").Value = -1293.7

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7142960.5
$ws.Range("I4").Value = 7142960.5
$ws.Range("K4").Value = 21428881.5
$ws.Range("M4").Value = -21428769.5
$ws.Range("H131").Value = 2863.16
$ws.Range("J131").Value = 2863.16
$ws.Range("L131").Value = 8589.48
$ws.Range("N131").Value = -18669.48

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 1054.2
$ws.Range("J19").Value = 1316.5
$ws.Range("L19").Value = 1316.5
$ws.Range("N19").Value = -1892.5
$ws.Range("H62").Value = 33705.168
$ws.Range("I62").Value = 37410.332
$ws.Range("K62").Value = 37410.332
$ws.Range("M62").Value = -36724.332
$ws.Range("H65").Value = 33705.168
$ws.Range("I65").Value = 37410.332
$ws.Range("K65").Value = 112230.996
$ws.Range("M65").Value = -108798.996
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H97").Value = 256.33334
$ws.Range("J97").Value = 195
$ws.Range("L97").Value = 195
$ws.Range("N97").Value = -1187
$ws.Range("H102").Value = 1125.1111
$ws.Range("I102").Value = 1166.375
$ws.Range("J102").Value = 795
$ws.Range("K102").Value = 1166.375
$ws.Range("L102").Value = 795
$ws.Range("M102").Value = 455.625
$ws.Range("N102").Value = -4039
$ws.Range("H132").Value = 898.5
$ws.Range("I132").Value = 898.5
$ws.Range("K132").Value = 2695.5
$ws.Range("M132").Value = -165.5
$ws.Range("H135").Value = 205438
$ws.Range("J135").Value = 205438
$ws.Range("L135").Value = 205438
$ws.Range("N135").Value = -215578

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1830.1666
$ws.Range("I68").Value = 1191.6666
$ws.Range("J68").Value = 2468.6667
$ws.Range("K68").Value = 1191.6666
$ws.Range("L68").Value = 2468.6667
$ws.Range("M68").Value = -442.6666
$ws.Range("N68").Value = -3966.6667
$ws.Range("H71").Value = 1830.1666
$ws.Range("I71").Value = 1191.6666
$ws.Range("J71").Value = 2468.6667
$ws.Range("K71").Value = 5958.333000000001
$ws.Range("L71").Value = 12343.3335
$ws.Range("M71").Value = -2214.333000000001
$ws.Range("N71").Value = -19831.3335
$ws.Range("H82").Value = 1200.6
$ws.Range("I82").Value = 1273.4286
$ws.Range("J82").Value = 1030.6666
$ws.Range("K82").Value = 1273.4286
$ws.Range("L82").Value = 1030.6666
$ws.Range("M82").Value = -912.4286
$ws.Range("N82").Value = -1752.6666
$ws.Range("H85").Value = 1200.6
$ws.Range("I85").Value = 1273.4286
$ws.Range("J85").Value = 1030.6666
$ws.Range("K85").Value = 1273.4286
$ws.Range("L85").Value = 1030.6666
$ws.Range("M85").Value = -25.42859999999996
$ws.Range("N85").Value = -3526.6666
$ws.Range("H132").Value = 3075.6155
$ws.Range("I132").Value = 2547.5
$ws.Range("K132").Value = 7642.5
$ws.Range("M132").Value = -5112.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2253750
$ws.Range("I4").Value = 1338333.4
$ws.Range("J4").Value = 5000000
$ws.Range("K4").Value = 1338333.4
$ws.Range("L4").Value = 5000000
$ws.Range("M4").Value = -1338220.4
$ws.Range("N4").Value = -5000226
$ws.Range("H100").Value = 3292.818
$ws.Range("I100").Value = 1802.8889
$ws.Range("K100").Value = 3605.7778
$ws.Range("M100").Value = -3064.7778
$ws.Range("H122").Value = 9099.799999999999
$ws.Range("I122").Value = 998
$ws.Range("J122").Value = 11125.25
$ws.Range("K122").Value = 2994
$ws.Range("L122").Value = 33375.75
$ws.Range("M122").Value = -544
$ws.Range("N122").Value = -38275.75
$ws.Range("H132").Value = 2248.9
$ws.Range("J132").Value = 3315.6667
$ws.Range("L132").Value = 9947.000100000001
$ws.Range("N132").Value = -15007.0001

Write-Output "Applied all cell updates."